$wb = $excel.ActiveWorkbook

# Sheet "展览" — rows 2,3,5,6,7,8,9 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 685
$ws1.Range("F3").Value = 29
$ws1.Range("F5").Value = 2001
$ws1.Range("F6").Value = 45
$ws1.Range("F7").Value = 3352
$ws1.Range("F8").Value = 462
$ws1.Range("F9").Value = 809

# Sheet "全部类型" — rows 2,3,6,7,8,9,10 in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 685
$ws4.Range("F3").Value = 29
$ws4.Range("F6").Value = 2001
$ws4.Range("F7").Value = 45
$ws4.Range("F8").Value = 3352
$ws4.Range("F9").Value = 462
$ws4.Range("F10").Value = 809
